$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "*YnBFawGccpQ*") {
        $p.Range.Font.Color = 5287936
        $p.Range.Font.ColorIndex = 0
    }
}
